{"js": "// Demote every top-level \"Heading 1\" section heading to \"Heading 2\"\n// (exercise 4 instructions: Objectives, Improve navigation and readability,\n// Implement dynamic calculations, Create and reference publication-ready\n// tables, Customise figures, Code, Add references, Finalise your MS Word\n// report).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.style === \"Heading 1\") {\n    paragraph.style = \"Heading 2\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Demote every top-level \"Heading 1\" section heading to \"Heading 2\"\n# (exercise 4 instructions: Objectives, Improve navigation and readability,\n# Implement dynamic calculations, Create and reference publication-ready\n# tables, Customise figures, Code, Add references, Finalise your MS Word\n# report).\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Heading 1\") {\n        $p.Style = \"Heading 2\"\n    }\n}\n"}
